# Update cryptocurrency prices and volume percentages (D/E columns)
# to the values reported by the GitHub Actions refresh run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.740.04"
$ws.Range("E2").Value = "'  -0.88%  "
$ws.Range("D3").Value = "'2.322.16"
$ws.Range("E3").Value = "'  -0.29%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'301.64"
$ws.Range("E5").Value = "'  -0.89%  "
$ws.Range("D6").Value = "'93.50"
$ws.Range("E6").Value = "'  -4.49%  "
$ws.Range("D7").Value = "'0.502"
$ws.Range("E7").Value = "'  -0.75%  "
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("D9").Value = "'0.492"
$ws.Range("E9").Value = "'  -2.25%  "
$ws.Range("D10").Value = "'33.86"
$ws.Range("E10").Value = "'  -4.98%  "
$ws.Range("D11").Value = "'0.0781"
$ws.Range("E11").Value = "'  -2.47%  "
$ws.Range("D12").Value = "'18.67"
$ws.Range("E12").Value = "'  -4.09%  "
$ws.Range("E13").Value = "'  +1.65%  "
$ws.Range("D14").Value = "'6.68"
$ws.Range("E14").Value = "'  -3.37%  "
$ws.Range("D15").Value = "'2.686.66"
$ws.Range("E15").Value = "'  -0.20%  "
$ws.Range("D16").Value = "'2.334.86"
$ws.Range("E16").Value = "'  +0.90%  "
$ws.Range("D17").Value = "'0.786"
$ws.Range("E17").Value = "'  -0.13%  "
$ws.Range("D18").Value = "'42.687.65"
$ws.Range("E19").Value = "'  -5.25%  "
$ws.Range("E20").Value = "'  +1.13%  "
$ws.Range("D21").Value = "'0.0₃0884"
$ws.Range("E21").Value = "'  -2.00%  "
$ws.Range("D22").Value = "'67.80"
$ws.Range("E22").Value = "'  -0.34%  "
$ws.Range("D23").Value = "'235.13"
$ws.Range("E23").Value = "'  -1.00%  "
$ws.Range("E24").Value = "'  +0.21%  "
$ws.Range("E25").Value = "'  +0.02%  "
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "'  -1.40%  "
$ws.Range("D27").Value = "'24.46"
$ws.Range("E27").Value = "'  -1.99%  "
$ws.Range("D28").Value = "'2.05"
$ws.Range("E28").Value = "'  -1.02%  "
$ws.Range("D29").Value = "'9.09"
$ws.Range("E29").Value = "'  -0.55%  "
$ws.Range("D30").Value = "'31.28"
$ws.Range("E30").Value = "'  -5.74%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  -0.01%  "
$ws.Range("D32").Value = "'139.75"
$ws.Range("E32").Value = "'  -15.92%  "
$ws.Range("E33").Value = "'  -0.58%  "
$ws.Range("D34").Value = "'17.43"
$ws.Range("E34").Value = "'  -3.11%  "
$ws.Range("D35").Value = "'0.0697"
$ws.Range("E35").Value = "'  +0.10%  "
$ws.Range("E36").Value = "'  -0.77%  "
$ws.Range("D37").Value = "'4.34"
$ws.Range("E37").Value = "'  -4.44%  "
$ws.Range("D38").Value = "'1.80"
$ws.Range("E38").Value = "'  +2.14%  "
$ws.Range("E39").Value = "'  -0.79%  "
$ws.Range("D40").Value = "'22.29"
$ws.Range("E40").Value = "'  +22.30%  "
$ws.Range("D41").Value = "'2.74"
$ws.Range("E41").Value = "'  -2.36%  "
$ws.Range("E42").Value = "'  -1.31%  "
$ws.Range("D43").Value = "'1.931.49"
$ws.Range("E43").Value = "'  -3.30%  "
$ws.Range("E44").Value = "'  -0.94%  "
$ws.Range("E45").Value = "'  -4.71%  "
$ws.Range("E46").Value = "'  -1.07%  "
$ws.Range("D47").Value = "'2.69"
$ws.Range("E47").Value = "'  -3.36%  "
$ws.Range("E48").Value = "'  +0.36%  "
$ws.Range("D49").Value = "'2.553.75"
$ws.Range("E49").Value = "'  -0.17%  "
$ws.Range("D50").Value = "'52.56"
$ws.Range("E50").Value = "'  -2.25%  "
$ws.Range("D51").Value = "'72.02"
$ws.Range("E51").Value = "'  +0.05%  "
